$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 473 - continuation of 9/13/2016 (serial 42626) entries
$ws.Range("A473").Value = "Demo"
$ws.Range("B473").Value = 42626
$ws.Range("C473").Value = "1830"
$ws.Range("D473").Value = "ACE"
$ws.Range("E473").Value = "013"

# Row 474 - new Skype Setup task
$ws.Range("A474").Value = "Skype Setup"
$ws.Range("B474").Value = 42626
$ws.Range("C474").Value = "1630"
$ws.Range("D474").Value = "OSG"
$ws.Range("E474").Value = "2028"
$ws.Range("F474").Value = "Skype. skype ID of the other party mccarthy-tetrault  x55401 Manusha. Login with itcyorku2, password York5065"
$ws.Rows.Item(474).RowHeight = 30

# Row 475 - new Skype Pickup task
$ws.Range("A475").Value = "Skype Pickup"
$ws.Range("B475").Value = 42626
$ws.Range("C475").Value = "1900"
$ws.Range("D475").Value = "OSG"
$ws.Range("E475").Value = "2028"
$ws.Range("F475").Value = "pick up skype kit, return to OSG 1014L"

# Rows 476-478 intentionally left blank (matches existing sheet's blank-row separators)

# Row 479 - 9/14/2016 (serial 42627) entries begin
$ws.Range("A479").Value = "Other"
$ws.Range("B479").Value = 42627
$ws.Range("C479").Value = "1730"
$ws.Range("D479").Value = "KT"
$ws.Range("E479").Value = "749"
$ws.Range("F479").Value = "Pick up portable DVD with cables and return to TEL 0003"

# Row 480
$ws.Range("A480").Value = "Demo"
$ws.Range("B480").Value = 42627
$ws.Range("C480").Value = "1845"
$ws.Range("D480").Value = "HNE"
$ws.Range("E480").Value = "104"

# Row 481
$ws.Range("A481").Value = "Demo"
$ws.Range("B481").Value = 42627
$ws.Range("C481").Value = "1900"
$ws.Range("D481").Value = "DB"
$ws.Range("E481").Value = "0004"

# Row 482
$ws.Range("A482").Value = "Demo"
$ws.Range("B482").Value = 42627
$ws.Range("C482").Value = "1900"
$ws.Range("D482").Value = "HNE"
$ws.Range("E482").Value = "401"

# Row 483 - new task type "2025"
$ws.Range("A483").Value = "Demo"
$ws.Range("B483").Value = 42627
$ws.Range("C483").Value = "2025"
$ws.Range("D483").Value = "OSG"
$ws.Range("E483").Value = "2008"

# Row 484
$ws.Range("A484").Value = "Demo"
$ws.Range("B484").Value = 42627
$ws.Range("C484").Value = "1630"
$ws.Range("D484").Value = "SSB"
$ws.Range("E484").Value = "W141"

# Row 485
$ws.Range("A485").Value = "Setup Mic"
$ws.Range("B485").Value = 42627
$ws.Range("C485").Value = "1630"
$ws.Range("D485").Value = "SSB"
$ws.Range("E485").Value = "W141"
$ws.Range("F485").Value = "2 neck mics and 2 hand held wireless mics for the pannel (there -check / test); podium mic (there); 2 wired audience mics (there /test). Please provide spare AA bateries for hand held"
$ws.Rows.Item(485).RowHeight = 45

# Row 486
$ws.Range("A486").Value = "Operator"
$ws.Range("B486").Value = 42627
$ws.Range("C486").Value = "1700"
$ws.Range("D486").Value = "SSB"
$ws.Range("E486").Value = "W141"
$ws.Range("F486").Value = "Operate event between 17:00-18:00"

# Row 487
$ws.Range("A487").Value = "AV Shutdown"
$ws.Range("B487").Value = 42627
$ws.Range("C487").Value = "2000"
$ws.Range("D487").Value = "SSB"
$ws.Range("E487").Value = "W141"

# Row 488
$ws.Range("A488").Value = "Pickup Mic"
$ws.Range("B488").Value = 42627
$ws.Range("C488").Value = "2000"
$ws.Range("D488").Value = "SSB"
$ws.Range("E488").Value = "W141"
$ws.Range("F488").Value = "return 2 wired audience mics, 2 hand held mics to rear booth"

# Move selection to the newly added last cell, matching the sheet's tracked state
$ws.Range("F488").Select()
